$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.833.40'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.740.03'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5172'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.88%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2810'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.23'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06117'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.754.44'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.50%  '
$ws.Range("E12").Value = '  +1.37%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.34'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6492'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.530'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '77.09'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9995'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '25.819.16'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.49'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006606'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.976.47'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.136'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.643'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.86%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.149'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '139.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.514'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.56%  '
$ws.Range("E28").Value = '  +0.64%  '
$ws.Range("E29").Value = '  +2.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '102.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08304'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.682'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.18%  '
$ws.Range("E33").Value = '  +1.65%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04502'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.611'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9881'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6162'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.655'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01584'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.941'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9990'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '100.79'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.3850'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7267'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.44%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.976'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.72%  '
$ws.Range("E46").Value = '  -1.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.296'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.19'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.690'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '29.92'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.22%  '
